$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Update the FilesTab Cypher query in cell B4: remove the File Type and Breed
# columns from the RETURN clause (ICDC Breed 1-14 script correction).
$newQuery = "MATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`nWHERE demo.breed IN ['Australian Shepherd']`nOPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`nOPTIONAL MATCH (samp:sample)-->(c)`nWITH DISTINCT f, parent, c, demo, diag, s`nRETURN  coalesce(f.file_name, '') AS ``File Name``,`n        coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(f.file_format, '') AS ``Format``,`n        coalesce(f.file_size, '') AS ``Size``,`n        coalesce(c.case_id, '') AS ``Case ID``,`n        coalesce(diag.disease_term,'') AS Diagnosis , `n        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

$ws.Range("B4").Value = $newQuery

# The removed lines shrink the wrapped text in B4, so the row shrinks to fit
# (matches Excel's automatic row-height recalculation for wrapped rows).
$ws.Rows.Item(4).RowHeight = 217.5

# Update the active selection on the sheet from D12 to D4.
$ws.Range("D4").Select()
